# BAARD uses fixed training set
# Update a handful of accuracy/FPR values on the "banknote" sheet to reflect
# the fixed training set results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("banknote")

$ws.Range("F4").Value = 78.57142857142857
$ws.Range("C5").Value = 94.3877551020408
$ws.Range("D12").Value = 91.83673469387756
$ws.Range("C13").Value = 79.08163265306123
$ws.Range("D14").Value = 56.63265306122449
$ws.Range("G19").Value = 2.040816326530612
